$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition list)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value  = 16778
$ws1.Range("F8").Value  = 124
$ws1.Range("F24").Value = 30
$ws1.Range("F26").Value = 6855
$ws1.Range("F36").Value = 4894

# Sheet "全部类型" (all types, combined list) - row numbers shifted by one vs above
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value  = 16778
$ws4.Range("F8").Value  = 124
$ws4.Range("F25").Value = 30
$ws4.Range("F27").Value = 6855
$ws4.Range("F38").Value = 4894
